$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

for ($r = 2; $r -le 73; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    $d = $epoch.AddDays($old)
    $next = $d.AddMonths(1)
    $newdate = Get-Date -Year $next.Year -Month $next.Month -Day 15 -Hour 0 -Minute 0 -Second 0
    $cell.Value2 = $newdate.ToOADate()
}
